$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prime the shared-string table so brand-new names land in the same order
# the original author's run produced them (first-seen-wins ordering).
$ws.Range("C7").Value = "T-Whisk"
$ws.Range("B4").Value = "Surf Co"

# Full data grid values (A1:E14) after the "create grid algo" shuffle.
$values = @(
    @("Monday", "Tuesday", "Wednesday", "Thursday", "Friday"),
    @("AM'er", "AM'er", "AM'er", "AM'er", "AM'er"),
    @("Smiles", "Jaws", "Nodder", "Jaws", "Chicken"),
    @("Nodder", "Surf Co", "Blister", "Nodder", "Jaws"),
    @("Aftie", "Aftie", "Aftie", "Aftie", "Aftie"),
    @("Surf Co", "Nodder", "Captain", "Blister", "Nodder"),
    @("Jaws", "Blister", "T-Whisk", "Captain", "Blister"),
    @("Wickie", "Wickie", "Wickie", "Wickie", "Wickie"),
    @("Blister", "Chicken", "Blister", "T-Whisk", "Captain"),
    @("Captain", "Captain", "Smiles", "Blister", "T-Whisk"),
    @("T-Whisk", "T-Whisk", "Surf Co", "Smiles", "Blister"),
    @("Blister", "O'fer", "O'fer", "O'fer", "Smiles"),
    @("O'fer", "Blister", "Chicken", "Surf Co", "O'fer"),
    @("Chicken", "Smiles", "Jaws", "Chicken", "Surf Co")
)

$cols = @("A", "B", "C", "D", "E")

for ($r = 0; $r -lt $values.Length; $r++) {
    $rowNum = $r + 1
    for ($c = 0; $c -lt 5; $c++) {
        $addr = "$($cols[$c])$rowNum"
        $ws.Range($addr).Value = $values[$r][$c]
    }
}

# Clear the two trailing rows that no longer exist in the shrunk grid.
$ws.Range("A15:E16").Clear()

# Bold header-style rows/cells (style index 1 in the XF table).
$boldCells = @("A1","B1","C1","D1","E1",
               "A2","B2","C2","D2","E2",
               "A5","B5","C5","D5","E5",
               "A8","B8","C8","D8","E8",
               "B12","C12","D12",
               "A13",
               "E13")
foreach ($addr in $boldCells) {
    $ws.Range($addr).Font.Bold = $true
}

# Small-font accent cells (style index 2 in the XF table).
$smallCells = @("C3","A4","B6","B9")
foreach ($addr in $smallCells) {
    $ws.Range($addr).Font.Size = 8
}

# Non-bold / non-small cells that previously carried a style and must revert to default.
$plainCells = @("A9","A14","B14","C14","D14")
foreach ($addr in $plainCells) {
    $ws.Range($addr).Font.Bold = $false
    $ws.Range($addr).Font.Size = 12
}

$ws.Range("F13").Select()
